# Add three new chart-data sheets (仪表盘/Gauge, 漏斗图/Funnel, 雷达图/Radar)
# to the workbook, matching the author's commit "add three graph(with bug)".

$wb = $excel.ActiveWorkbook

# --- 仪表盘 (Gauge) -----------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sGauge = $wb.Worksheets.Add($null, $lastSheet)
$sGauge.Name = "仪表盘"
$sGauge.Range("A1").Value = "完成率"
$sGauge.Range("B1").Value = 55.5

# --- 漏斗图 (Funnel) -----------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sFunnel = $wb.Worksheets.Add($null, $lastSheet)
$sFunnel.Name = "漏斗图"
$sFunnel.Range("A1").Value = "展现"
$sFunnel.Range("B1").Value = 100
$sFunnel.Range("A2").Value = "点击"
$sFunnel.Range("B2").Value = 80
$sFunnel.Range("A3").Value = "访问"
$sFunnel.Range("B3").Value = 60
$sFunnel.Range("A4").Value = "咨询"
$sFunnel.Range("B4").Value = 40
$sFunnel.Range("A5").Value = "订单"
$sFunnel.Range("B5").Value = 20

# --- 雷达图 (Radar) -------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sRadar = $wb.Worksheets.Add($null, $lastSheet)
$sRadar.Name = "雷达图"

# Row labels + first two data columns entered first (matches the shared-
# string ordering produced by the original authoring session), header
# row for columns C/D added afterwards.
$sRadar.Range("A2").Value = "销售"
$sRadar.Range("B2").Value = 6500
$sRadar.Range("C2").Value = 4300
$sRadar.Range("D2").Value = 5000

$sRadar.Range("A3").Value = "经理"
$sRadar.Range("B3").Value = 16000
$sRadar.Range("C3").Value = 10000
$sRadar.Range("D3").Value = 14000

$sRadar.Range("A4").Value = "信息技术"
$sRadar.Range("B4").Value = 30000
$sRadar.Range("C4").Value = 28000
$sRadar.Range("D4").Value = 28000

$sRadar.Range("A5").Value = "客服"
$sRadar.Range("B5").Value = 38000
$sRadar.Range("C5").Value = 35000
$sRadar.Range("D5").Value = 31000

$sRadar.Range("A6").Value = "研发"
$sRadar.Range("B6").Value = 52000
$sRadar.Range("C6").Value = 50000
$sRadar.Range("D6").Value = 42000

$sRadar.Range("A7").Value = "市场"
$sRadar.Range("B7").Value = 25000
$sRadar.Range("C7").Value = 19000
$sRadar.Range("D7").Value = 21000

# Header row for the radar chart's two value series (column B is left
# without a header - this is the "bug" referenced in the commit message).
$sRadar.Range("C1").Value = "预算分配"
$sRadar.Range("D1").Value = "实际开销"
